$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet index 1) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 854
$ws1.Range("F6").Value = 14
$ws1.Range("F7").Value = 3852
$ws1.Range("F8").Value = 2548
$ws1.Range("F10").Value = 2408
$ws1.Range("F14").Value = 1629
$ws1.Range("F16").Value = 6
$ws1.Range("F21").Value = 265
$ws1.Range("F23").Value = 438
$ws1.Range("F24").Value = 26
$ws1.Range("F26").Value = 486
$ws1.Range("F27").Value = 669
$ws1.Range("F30").Value = 370
$ws1.Range("F33").Value = 875
$ws1.Range("F34").Value = 41
$ws1.Range("F36").Value = 934
$ws1.Range("F37").Value = 1957
$ws1.Range("F39").Value = 518
$ws1.Range("F40").Value = 79
$ws1.Range("F42").Value = 595
$ws1.Range("F43").Value = 1241
$ws1.Range("F44").Value = 33
$ws1.Range("F46").Value = 413

# --- Sheet "全部类型" (sheet index 4) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 854
$ws4.Range("F5").Value = 14
$ws4.Range("F6").Value = 3852
$ws4.Range("F7").Value = 2548
$ws4.Range("F8").Value = 2408
$ws4.Range("F10").Value = 1629
$ws4.Range("F13").Value = 6
$ws4.Range("F18").Value = 265
$ws4.Range("F20").Value = 438
$ws4.Range("F21").Value = 26
$ws4.Range("F23").Value = 486
$ws4.Range("F24").Value = 669
$ws4.Range("F30").Value = 370
$ws4.Range("C31").Value = "杭州·第37届 中二病 原神x星穹only"
$ws4.Range("D31").Value = "康候圣街99号 顺丰创新中心"
$ws4.Range("E31").Value = "2024.05.01 10:30-05.02 17:00"
$ws4.Range("F31").Value = 1609
$ws4.Range("G31").Value = 60
$ws4.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=82700"
$ws4.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202403/Kb75MESZ1710215541381.jpeg"
$ws4.Range("C32").Value = "杭州·第7届YH樱花动漫游戏文化节"
$ws4.Range("D32").Value = "德胜东路2539号 梦马汽车小镇"
$ws4.Range("E32").Value = "2024.05.01 10:00-05.02 17:00"
$ws4.Range("F32").Value = 875
$ws4.Range("G32").Value = 70
$ws4.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=82828"
$ws4.Range("I32").Value = "//i1.hdslb.com/bfs/openplatform/202403/Kd0niodt1710905544733.jpeg"
$ws4.Range("C33").Value = "杭州·第7届YH樱花漫展-SVIP嘉宾前排票"
$ws4.Range("F33").Value = 41
$ws4.Range("G33").Value = 168
$ws4.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=83267"
$ws4.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202403/DgmIZ6G71711357279757.jpeg"
$ws4.Range("C34").Value = "杭州·第7届YH樱花漫展-配音演员紫枫儿内场票"
$ws4.Range("E34").Value = "2024.05.01 10:00-05.01 17:00"
$ws4.Range("F34").Value = 8
$ws4.Range("G34").Value = 98
$ws4.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=83331"
$ws4.Range("I34").Value = "//i0.hdslb.com/bfs/openplatform/202403/h5ilz3SA1711351453471.jpeg"
$ws4.Range("F36").Value = 934
$ws4.Range("F37").Value = 1957
$ws4.Range("F42").Value = 518
$ws4.Range("F43").Value = 79
$ws4.Range("F45").Value = 595
$ws4.Range("F46").Value = 1241
$ws4.Range("F47").Value = 33
$ws4.Range("F48").Value = 413
